$d = $word.ActiveDocument

# 1. "Joosd during his lunch break..." paragraph: the old text was split across
#    three runs with grammar-check proofErr markers around "is". Collapse it
#    back down to a single clean run with the same final wording.
$d.Content.Find.Execute(
    "Joosd during his lunch break wants to check if there is any upcoming concerts.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Joosd during his lunch break wants to check if there is any upcoming concerts.",
    2)

# 2. Append a trailing comma after each of the colour swatches and after the
#    font-family line (so the list reads as a comma separated list).
$d.Content.Find.Execute("#3c415e", $true, $false, $false, $false, $false, $true, 1, $false, "#3c415e,", 2)
$d.Content.Find.Execute("#738598", $true, $false, $false, $false, $false, $true, 1, $false, "#738598,", 2)
$d.Content.Find.Execute("#dfe2e2", $true, $false, $false, $false, $false, $true, 1, $false, "#dfe2e2,", 2)
$d.Content.Find.Execute("#1cb3c8", $true, $false, $false, $false, $false, $true, 1, $false, "#1cb3c8,", 2)
$d.Content.Find.Execute("font-family: Agency FB", $true, $false, $false, $false, $false, $true, 1, $false, "font-family: Agency FB,", 2)

# 3. The empty paragraph right after the font-family line now documents the
#    font-size deliverable too. Clone the run formatting (lang=en-US) from
#    the font-family paragraph itself so the new run carries the same
#    <w:rPr/>, then overwrite its text with "font-size: 30px".
$i = 1
$sourceIndex = -1
while ($i -le $d.Paragraphs.Count) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text.TrimEnd() -eq "font-family: Agency FB,") {
        $sourceIndex = $i
        $i = $d.Paragraphs.Count + 1
    } else {
        $i = $i + 1
    }
}
$targetIndex = $sourceIndex + 1
$formatted = $d.Paragraphs.Item($sourceIndex).Range.FormattedText
$target = $d.Paragraphs.Item($targetIndex)
$target.Range.FormattedText = $formatted
$target = $d.Paragraphs.Item($targetIndex)
$textOnly = $d.Range($target.Range.Start, $target.Range.End - 1)
$textOnly.Text = "font-size: 30px"
